# mactaquac-egg_development.xlsx edit:
#   - Picking sheet: add "First Hatch Observed (Y/N)" / "100% Hatch Observed (Y/N)"
#     columns, replace the numbered pick-day columns with two example date
#     headers + helper cells showing the expected text/date formats, update
#     the explanatory cell comments, and make "Picking" the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Picking")

# ---------------------------------------------------------------------
# 1. Make room for the two new "Hatch Observed" columns by inserting two
#    blank columns in front of the old "Comments" column (K). Everything
#    that used to live in K:P (Comments, 1, 2, 3, 4, 5) shifts right to
#    M:R, carrying its formatting with it.
# ---------------------------------------------------------------------
$ws.Range("K1:L1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2. The "pick day" columns used to just be numbered 1-5; now only two
#    example date headers remain, showing the expected date format.
# ---------------------------------------------------------------------
$ws.Range("N3").Value = "1999-Jan-1"
$ws.Range("O3").Value = "1999-Jan-2"

# Drop the old trailing numbered headers (now shifted to P3:R3).
$ws.Range("P3:R3").Clear()

# ---------------------------------------------------------------------
# 3. New headers for the two inserted columns.
# ---------------------------------------------------------------------
$ws.Range("K3").Value = "First Hatch Observed (Y/N)"
$ws.Range("L3").Value = "100% Hatch Observed (Y/N)"

# ---------------------------------------------------------------------
# 4. Helper cells below the new date-example headers, showing the
#    expected data format for that column: an actual date entry first,
#    then a free-text entry further down (order matters for how the
#    shared style table gets built up).
# ---------------------------------------------------------------------
$ws.Range("K6").NumberFormat = "d-mmm-yy"
$ws.Range("K4").NumberFormat = "@"

# ---------------------------------------------------------------------
# 5. Update the cell comments: the old single "Day picks were made"
#    comment on the (now repurposed) L3 cell is replaced by matching
#    comments on the new N3/O3 date headers.
# ---------------------------------------------------------------------
$ws.Range("L3").Comment.Delete()
$ws.Range("N3").AddComment("Date picks were made.`nFormat: YYYY-mmm-dd")
$ws.Range("O3").AddComment("Date picks were made.`nFormat: YYYY-mmm-dd")

# ---------------------------------------------------------------------
# 6. Resize the affected columns to fit their new headers.
# ---------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 18.5703125
$ws.Columns.Item(9).ColumnWidth = 14.140625
$ws.Columns.Item(10).ColumnWidth = 10.28515625
$ws.Columns.Item(11).ColumnWidth = 24.85546875
$ws.Columns.Item(12).ColumnWidth = 25.7109375
$ws.Columns.Item(13).ColumnWidth = 10.5703125
$ws.Columns.Item(14).ColumnWidth = 15
$ws.Columns.Item(15).ColumnWidth = 10.28515625

# ---------------------------------------------------------------------
# 7. "Picking" becomes the active/visible sheet, with L3 selected.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("L3").Select()
